# The experiment was regenerated with updated distance/size values.
# Distances: D80 -> D86, D64 -> D69, D51 -> D55
# Sizes:     S30 -> S31
# These substrings appear throughout the Condition / Filename_Left /
# Filename_Right / Distance / Size columns (all stored as shared strings),
# so a workbook-wide text replace reproduces the regenerated order file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("S30", "S31")
